$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update the "ambiente" (environment) row 2: old ssurgwsoadev4 domain -> ssurgwsoadev4-oci ---
# Use a scratch cell to preserve A2's existing cell format (it uses the "quote prefix" style)
# across the value change, since re-assigning .Value normally resets formatting.
$ws.Range("ZZ1").Value = "x"
$ws.Range("A2").Copy() | Out-Null
$ws.Range("ZZ1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"

$ws.Range("ZZ1").Copy() | Out-Null
$ws.Range("A2").PasteSpecial(-4122) | Out-Null    # xlPasteFormats restored on A2
$ws.Range("ZZ1").Clear() | Out-Null
$excel.CutCopyMode = $false

# B2 holds the matching URL text + a hyperlink to the same address
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"

foreach ($hl in $ws.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.Address = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
    }
}

# --- 2) Update the numeric "Documento" values ---
$ws.Range("G2").Value = 30694790255
$ws.Range("G3").Value = 30522093323
$ws.Range("M3").Value = 303

# --- 3) Update the view: selection moves to M4, window scrolled so column D is leftmost ---
$ws.Range("M4").Select()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1

Write-Output "edits applied"
